# Analysis as of Sept 2021 Submission
# Update the "Demographics" sheet (race/ethnicity rows and mPSG severity rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# --- Race/ethnicity block (rows 12-14): cyclic re-ordering of rows ---
# Row 12 "native american" -> becomes "pacific islander" (values unchanged)
$ws.Range("A12").Value = "pacific islander"
$ws.Range("B12").Value = "3/510 (0.6%)"
$ws.Range("C12").Value = "0/11 (0.0%)"
$ws.Range("D12").Value = "1/62 (1.6%)"
$ws.Range("E12").Value = "1/266 (0.4%)"
$ws.Range("F12").Value = "1/171 (0.6%)"

# Row 13 "asian" -> becomes "native american" (D/F values updated)
$ws.Range("A13").Value = "native american"
$ws.Range("B13").Value = "3/510 (0.6%)"
$ws.Range("C13").Value = "0/11 (0.0%)"
$ws.Range("D13").Value = "1/62 (1.6%)"
$ws.Range("E13").Value = "1/266 (0.4%)"
$ws.Range("F13").Value = "1/171 (0.6%)"

# Row 14 "pacific islander" -> becomes "asian" (D/F values updated)
$ws.Range("A14").Value = "asian"
$ws.Range("B14").Value = "3/510 (0.6%)"
$ws.Range("C14").Value = "0/11 (0.0%)"
$ws.Range("D14").Value = "0/62 (0.0%)"
$ws.Range("E14").Value = "1/266 (0.4%)"
$ws.Range("F14").Value = "2/171 (1.2%)"

# --- mPSG severity block (rows 35-36): updated counts/percentages ---
# Row 35 "mild"
$ws.Range("B35").Value = "93/510 (18.2%)"
$ws.Range("E35").Value = "55/266 (20.7%)"
$ws.Range("F35").Value = "30/171 (17.5%)"

# Row 36 "none"
$ws.Range("B36").Value = "0/510 (0.0%)"
$ws.Range("E36").Value = "0/266 (0.0%)"
$ws.Range("F36").Value = "0/171 (0.0%)"
